# Generate Report for Handoff
# Adds a new row (row 3) to each of the three tables (Overview, zh-cn, de-de)
# describing the hand-off of a new file:
#   58fdcd79-eaf1-4e54-b01b-339e81494d1eo...(padding)...md

$wb = $excel.ActiveWorkbook

$mdName   = '58fdcd79-eaf1-4e54-b01b-339e81494d1eooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$mdPath   = 'e2e\58fdcd79-eaf1-4e54-b01b-339e81494d1eooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$mdUrl    = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cb76949b6819939256c6d52f1a1494b8f633ebc3/e2e/58fdcd79-eaf1-4e54-b01b-339e81494d1eooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'

$zhXlf    = '58fdcd79-eaf1-4e54-b01b-339e81494d1eoooooooooooooooooooooooooooooooooooooooo.7ce1c7bf275b65b7161cc1e3aa25392a1541e048.zh-cn.xlf'
$deXlf    = '58fdcd79-eaf1-4e54-b01b-339e81494d1eoooooooooooooooooooooooooooooooooooooooo.7ce1c7bf275b65b7161cc1e3aa25392a1541e048.de-de.xlf'

$status   = 'Ready for handoff'
$extension = '.md'

$dtOverviewDe = '2016-08-16 20:24:50'
$dtZh         = '2016-08-16 20:24:45'
$dtDe         = '2016-08-16 20:24:50'

$hyperLinkColor = 15570276  # BGR for FF6495ED, matching the existing HyperLink style

function Style-Hyperlink($range) {
    $range.Style = "HyperLink"
    $range.Font.Underline = $true
    $range.Font.Color = $hyperLinkColor
}

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) -> table "Overview"
# Columns: A File Name | B Path And Name | C Extension | D Publish URL |
#          E zh-cn | F de-de | G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $mdName
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $mdUrl, "", "", $mdPath) | Out-Null
Style-Hyperlink $wsOverview.Range("B3")
$wsOverview.Range("C3").Value = $extension
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("G3").Value = $dtOverviewDe

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) -> table "zh-cn"
# Columns: A Source File Name | B File Extension | C Status | D Source Path |
#          E Priority | F Content Duplicate | G Latest Handoff File |
#          H Latest Handoff Datetime | I Latest Target File |
#          J Latest Handback File | K Latest Handback DateTime |
#          L Reference Tokens | M To be localized | N Dependency From |
#          O Has metadata | P Error Detail
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl, "", "", $mdName) | Out-Null
Style-Hyperlink $wsZh.Range("A3")
$wsZh.Range("B3").Value = $extension
$wsZh.Range("C3").Value = $status
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H3").Value = $dtZh
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3) -> table "de-de"
# Same column layout as zh-cn
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl, "", "", $mdName) | Out-Null
Style-Hyperlink $wsDe.Range("A3")
$wsDe.Range("B3").Value = $extension
$wsDe.Range("C3").Value = $status
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H3").Value = $dtDe
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

# Widen columns E/F (zh-cn/de-de date columns) on Overview and column C
# (Status) on the zh-cn / de-de sheets to fit the new content, matching the
# author's column-width tweak (old width ~13.41 -> new width ~17.22 chars).
$wsOverview.Columns("E").ColumnWidth = 17.1
$wsOverview.Columns("F").ColumnWidth = 17.1
$wsZh.Columns("C").ColumnWidth = 17.1
$wsDe.Columns("C").ColumnWidth = 17.1
